$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update reference list to include R15, and bump quantity from 2 to 3
$ws.Range("C17").Value = "R6, R13, R15"
$ws.Range("B17").Value = 3

# Update the active cell selection as the user would have left it
$ws.Range("C18").Select()
